$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.447.25'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.882.21'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'0.7185"
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = "'243.92"
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = "'0.07946"
$ws.Range('E8').Value = '  +2.01%  '
$ws.Range('D9').Value = "'0.3156"
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').Value = "'25.03"
$ws.Range('D11').Value = "'0.08143"
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('D12').Value = '1.897.21'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').Value = "'5.262"
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = "'94.81"
$ws.Range('E14').Value = '  +4.12%  '
$ws.Range('D15').Value = "'0.7110"
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = "'6.411"
$ws.Range('E16').Value = '  +4.46%  '
$ws.Range('D17').Value = "'0.000008443"
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').Value = '29.471.95'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = "'254.34"
$ws.Range('E19').Value = '  +5.94%  '
$ws.Range('D20').Value = "'13.34"
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').Value = '2.142.58'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = "'7.777"
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').Value = "'1.001"
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = "'0.1588"
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = "'9.092"
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').Value = "'162.97"
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = "'18.96"
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').Value = "'4.429"
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = "'1.229"
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = "'0.05342"
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = "'1.956"
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').Value = "'0.7574"
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('D37').Value = "'2.703"
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').Value = "'0.01901"
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('D39').Value = '1.276.15'
$ws.Range('E39').Value = '  +2.57%  '
$ws.Range('D40').Value = "'2.771"
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').Value = "'6.471"
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('D42').Value = "'112.99"
$ws.Range('E42').Value = '  +3.34%  '
$ws.Range('D43').Value = "'74.49"
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('D44').Value = "'0.9048"
$ws.Range('E44').Value = '  +1.41%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = "'1.002"
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000130"
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('D47').Value = '2.037.79'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').Value = "'0.5208"
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = "'9.502"
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').Value = "'0.4379"
$ws.Range('E51').Value = '  +0.84%  '

# Reset style on cells where a leading apostrophe was used to force text,
# so no stray quotePrefix style gets attached to the cell.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
